$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the green fill (FF92D050) to the C4:E32 homework-score block.
# This both creates the new fill/style entries and re-stamps every cell
# in the range with the new style id (matching s="2" -> s="6").
$ws.Range("C4:E32").Interior.Color = 5296274

# Cells that were previously blank in C/D/E now carry an explicit 0.
$zeroCells = @(
    "C4", "E4",
    "C8",
    "C14", "D14", "E14",
    "C15", "D15", "E15",
    "E21",
    "C24", "D24", "E24",
    "E27",
    "C30"
)
foreach ($addr in $zeroCells) {
    $ws.Range($addr).Value = 0
}

# New column L ("Лаб_1" extra makeup work?) with explicit 0 entries for the
# same four rows that previously had no data at all (11, 21, 22, 24 -> rows 14,15,24 and 4).
$lCells = @("L4", "L14", "L15", "L24")
foreach ($addr in $lCells) {
    $ws.Range($addr).Value = 0
}

# Restore the final selection to match the authored view state.
$ws.Range("L8").Select()
